$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.552.63"

$ws.Range("D3").Value = "2.069.27"
$ws.Range("E3").Value = "  -0.38%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.39"
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.48"
$ws.Range("E8").Value = "  -2.16%  "

$ws.Range("E9").Value = "  -1.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0779"
$ws.Range("E10").Value = "  -0.60%  "

$ws.Range("E11").Value = "  +1.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.87"
$ws.Range("E12").Value = "  +0.39%  "

$ws.Range("D13").Value = "2.374.86"
$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.89"
$ws.Range("E14").Value = "  -0.41%  "

$ws.Range("E15").Value = "  -1.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.32"
$ws.Range("E16").Value = "  -0.96%  "

$ws.Range("D17").Value = "2.061.30"
$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("D18").Value = "37.503.10"
$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.51"
$ws.Range("E19").Value = "  -1.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.93"
$ws.Range("E20").Value = "  -3.10%  "

$ws.Range("E21").Value = "  -0.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.04"
$ws.Range("E22").Value = "  -0.19%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  -1.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.61"
$ws.Range("E26").Value = "  +5.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.27"
$ws.Range("E27").Value = "  -1.24%  "

$ws.Range("E28").Value = "  -3.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.40"
$ws.Range("E29").Value = "  -0.57%  "

$ws.Range("E30").Value = "  -1.88%  "

$ws.Range("E31").Value = "  +1.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.59"
$ws.Range("E32").Value = "  -1.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0632"
$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.60"
$ws.Range("E34").Value = "  -1.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.46"
$ws.Range("E35").Value = "  -1.27%  "

$ws.Range("E36").Value = "  -0.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.31"
$ws.Range("E37").Value = "  -2.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.28"
$ws.Range("E39").Value = "  -1.42%  "

$ws.Range("E40").Value = "  +6.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.42"
$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("E42").Value = "  +4.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0956"

$ws.Range("E44").Value = "  +0.84%  "

$ws.Range("D45").Value = "1.471.97"
$ws.Range("E45").Value = "  +2.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.60"
$ws.Range("E46").Value = "  -0.31%  "

$ws.Range("E47").Value = "  -1.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.98"
$ws.Range("E48").Value = "  -5.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.24"
$ws.Range("E49").Value = "  -2.02%  "

$ws.Range("E50").Value = "  -1.87%  "

$ws.Range("D51").Value = "2.259.34"
$ws.Range("E51").Value = "  -0.32%  "
